$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 15 de Septiembre de 2020 a las 16:46"

# Full target table (country + 7 numeric stat columns) for rows 4..219,
# reflecting both the reordering of Portugal/Zambia/Uganda/Birmania and
# the refreshed case counts for all affected countries.
$data = @(
    @('Estados Unidos',6751889,2600,4029052,2523673,0,164,199164),
    @('India',4963097,36183,3887371,994558,0,360,81168),
    @('Brasil',4349544,0,3613184,604243,0,0,132117),
    @('Rusia',1073849,5529,884305,170759,0,150,18785),
    @('Peru',733860,0,573364,129684,0,0,30812),
    @('Colombia',721892,0,606925,91844,0,0,23123),
    @('Mexico',671716,3335,475795,124872,0,228,71049),
    @('Sudafrica',650749,0,579289,55961,0,0,15499),
    @('España',593730,0,0,0,0,0,29848),
    @('Argentina',565446,0,438883,114853,0,43,11710),
    @('Chile',436433,0,407725,16695,0,0,12013),
    @('Iran',407353,2705,349984,33916,0,140,23453),
    @('Francia',387252,0,89507,266795,0,0,30950),
    @('Reino Unido',371125,0,0,0,0,0,41637),
    @('Banglades',341056,1724,245594,90660,0,43,4802),
    @('Arabia Saudita',326930,672,305022,17570,0,33,4338),
    @('Pakistan',302424,404,290261,5774,0,6,6389),
    @('Irak',294478,0,229132,57260,0,0,8086),
    @('Turquia',292878,0,260058,25701,0,0,7119),
    @('Italia',288761,0,213950,39187,0,0,35624),
    @('Filipinas',269407,3544,207352,57392,0,34,4663),
    @('Alemania',263954,733,237550,16967,0,1,9437),
    @('Indonesia',225030,3507,161065,55000,0,124,8965),
    @('Israel',162273,1905,120443,40689,0,5,1141),
    @('Ucrania',159702,2905,70810,85628,0,53,3264),
    @('Canada',138010,0,121224,7607,0,0,9179),
    @('Bolivia',127619,828,85198,35027,0,50,7394),
    @('Catar',122214,239,119144,2862,0,1,208),
    @('Ecuador',118911,0,97063,10926,0,0,10922),
    @('Kazajistan',106920,65,100836,4450,0,0,1634),
    @('Rumania',105298,1111,43244,57818,0,51,4236),
    @('Republica Dominicana',104110,0,77790,24336,0,0,1984),
    @('Panama',102204,0,74782,25249,0,0,2173),
    @('Egipto',101177,0,84969,10547,0,0,5661),
    @('Kuwait',96301,829,86219,9514,0,5,568),
    @('Belgica',94306,851,18737,65642,0,2,9927),
    @('Oman',90660,438,84113,5750,0,7,797),
    @('Marruecos',88203,0,68970,17619,0,0,1614),
    @('Suecia',87345,0,0,0,0,4,5851),
    @('China',85202,8,80426,142,0,0,4634),
    @('Paises Bajos',84778,1379,0,0,0,2,6258),
    @('Guatemala',82172,0,71352,7848,0,0,2972),
    @('Emiratos Arabes Unidos',80940,674,70635,9904,0,2,401),
    @('Japon',75657,0,67242,6973,0,0,1442),
    @('Polonia',75134,605,61548,11359,0,24,2227),
    @('Bielorrusia',74552,192,72661,1130,0,5,761),
    @('Honduras',68620,831,18487,48046,0,8,2087),
    @('Portugal',65021,425,44362,18784,0,4,1875),
    @('Etiopia',64786,0,25333,38431,0,0,1022),
    @('Venezuela',61569,0,49371,11704,0,0,494),
    @('Barein',60965,0,54204,6548,0,0,213),
    @('Singapur',57488,34,56802,659,0,0,27),
    @('Costa Rica',57361,0,21206,35534,0,0,621),
    @('Nepal',56788,1459,40638,15779,0,11,371),
    @('Nigeria',56388,0,44337,10968,0,0,1083),
    @('Argelia',48496,0,34204,12672,0,0,1620),
    @('Uzbekistan',48429,593,44942,3085,0,6,402),
    @('Suiza',47751,315,39600,6126,0,0,2025),
    @('Armenia',46119,150,41941,3258,0,1,920),
    @('Ghana',45601,0,44679,628,0,0,294),
    @('Kirguistan',44999,71,41103,2833,0,0,1063),
    @('Moldavia',43734,527,30437,12154,0,14,1143),
    @('Afganistan',38815,43,32098,5291,0,1,1426),
    @('Azerbaiyan',38517,114,35998,1953,0,2,566),
    @('Chequia',37222,0,22020,14737,0,0,465),
    @('Kenia',36205,0,23243,12338,0,0,624),
    @('Austria',34305,764,27354,6194,0,0,757),
    @('Serbia',32511,74,31313,463,0,2,735),
    @('Estado de Palestina',32250,888,21804,10217,0,3,229),
    @('Irlanda',31192,0,23364,6044,0,0,1784),
    @('Paraguay',28367,0,14814,13014,0,0,539),
    @('El Salvador',27088,79,19450,6846,0,4,792),
    @('Australia',26739,47,23652,2271,0,0,816),
    @('Libano',24857,0,8765,15846,0,0,246),
    @('Bosnia y Herzegovina',23929,294,16701,6503,0,20,725),
    @('Libia',23515,0,12762,10385,0,0,368),
    @('Corea del Sur',22391,106,18878,3146,0,4,367),
    @('Dinamarca',20571,334,16557,3381,0,0,633),
    @('Camerun',20228,0,18837,976,0,0,415),
    @('Costa de Marfil',19066,0,18174,772,0,0,120),
    @('Bulgaria',18061,0,12930,4402,0,0,729),
    @('Republica de Macedonia',15925,98,13418,1846,0,9,661),
    @('Madagascar',15803,34,14452,1137,0,1,214),
    @('Senegal',14529,223,10692,3539,0,1,298),
    @('Hungria',13879,726,4130,9103,0,4,646),
    @('Zambia',13819,99,12590,905,0,4,324),
    @('Croacia',13749,151,11412,2107,0,3,230),
    @('Sudan',13535,0,6759,5940,0,0,836),
    @('Grecia',13420,0,3804,9306,0,0,310),
    @('Noruega',12330,54,10371,1694,0,0,265),
    @('Albania',11520,0,6615,4567,0,0,338),
    @('Consejo Danes para los Refugiados',10401,11,9817,317,0,3,267),
    @('Guinea',10061,0,9352,646,0,0,63),
    @('Malasia',9969,23,9209,632,0,0,128),
    @('Namibia',9818,0,6693,3022,0,0,103),
    @('Guayana Francesa',9552,0,9156,333,0,0,63),
    @('Maldivas',9243,0,7536,1674,0,0,33),
    @('Tayikistan',9129,41,7896,1160,0,1,73),
    @('Finlandia',8725,98,7500,886,0,2,339),
    @('Gabon',8654,0,7785,816,0,0,53),
    @('Haiti',8499,0,6120,2160,0,0,219),
    @('Zimbabue',7531,0,5690,1617,0,0,224),
    @('Tunez',7382,0,2175,5090,0,0,117),
    @('Mauritania',7295,0,6835,299,0,0,161),
    @('Luxemburgo',7244,0,6555,565,0,0,124),
    @('Montenegro',6900,0,4589,2190,0,0,121),
    @('Eslovaquia',5768,188,3214,2516,0,0,38),
    @('Malaui',5697,0,3742,1777,0,0,178),
    @('Mozambique',5482,0,3024,2423,0,0,35),
    @('Republica de Yibuti',5396,0,5331,4,0,0,61),
    @('Uganda',5123,145,2333,2732,0,2,58),
    @('Suazilandia',5104,0,4374,629,0,0,101),
    @('Guinea Ecuatorial',5000,0,4496,421,0,0,83),
    @('Hong Kong',4976,4,4646,228,0,1,102),
    @('Congo',4934,0,3887,959,0,0,88),
    @('Cabo Verde',4839,0,4240,554,0,0,45),
    @('Nicaragua',4818,0,2913,1761,0,0,144),
    @('Republica de Africa Central',4772,0,1828,2882,0,0,62),
    @('Cuba',4726,0,4040,578,0,0,108),
    @('Surinam',4611,0,3935,581,0,0,95),
    @('Ruanda',4602,0,2736,1844,0,0,22),
    @('Jamaica',3933,0,1161,2728,0,0,44),
    @('Eslovenia',3831,82,2789,907,0,0,135),
    @('Siria',3576,0,858,2561,0,0,157),
    @('Jordania',3528,0,2255,1247,0,0,26),
    @('Birmania',3502,307,832,2635,0,3,35),
    @('Tailandia',3480,5,3315,107,0,0,58),
    @('Angola',3439,0,1324,1979,0,0,136),
    @('Gambia',3405,0,1723,1579,0,0,103),
    @('Lituania',3397,11,2094,1216,0,0,87),
    @('Somalia',3389,0,2803,488,0,0,98),
    @('Mayotte',3374,0,2964,370,0,0,40),
    @('Sri Lanka',3271,9,3016,242,0,0,13),
    @('Trinidad yTobago',3183,42,798,2329,0,1,56),
    @('Guadalupe',3080,0,837,2219,0,0,24),
    @('Aruba',3060,0,1566,1474,0,0,20),
    @('Bahamas',3008,34,1391,1549,0,1,68),
    @('Mali',2935,0,2289,518,0,0,128),
    @('Reunion',2872,0,1313,1544,0,0,15),
    @('Estonia',2722,25,2286,372,0,0,64),
    @('Sudan del Sur',2587,0,1290,1248,0,0,49),
    @('Georgia',2562,170,1369,1174,0,0,19),
    @('Botsuana',2463,0,575,1877,0,0,11),
    @('Malta',2454,49,1931,507,0,0,16),
    @('Guinea-Bisau',2275,0,1127,1109,0,0,39),
    @('Benin',2267,0,1942,285,0,0,40),
    @('Islandia',2174,6,2102,62,0,0,10),
    @('Sierra Leona',2111,0,1636,403,0,0,72),
    @('Yemen',2013,0,1215,215,0,0,583),
    @('Guyana',1884,0,1265,563,0,0,56),
    @('Uruguay',1812,0,1528,239,0,0,45),
    @('Nueva Zelanda',1801,3,1694,83,0,0,24),
    @('Burkina Faso',1717,0,1137,524,0,0,56),
    @('Togo',1578,0,1204,334,0,0,40),
    @('Republica de Chipre',1534,0,1282,230,0,0,22),
    @('Belice',1501,21,540,942,0,0,19),
    @('Letonia',1482,5,1248,199,0,0,35),
    @('Principado de Andorra',1438,0,945,440,0,0,53),
    @('Liberia',1321,0,1213,26,0,0,82),
    @('Lesoto',1245,0,528,684,0,0,33),
    @('Niger',1180,0,1104,7,0,0,69),
    @('Polinesia Francesa',1099,0,672,425,0,0,2),
    @('Republica del Chad',1085,0,940,64,0,0,81),
    @('Vietnam',1063,0,926,102,0,0,35),
    @('Martinica',939,0,98,823,0,0,18),
    @('Santo Tome y Principe',906,0,870,21,0,0,15),
    @('San Marino',722,0,662,18,0,0,42),
    @('Crucero',712,0,651,48,0,0,13),
    @('Islas Turcas y Caicos',648,2,532,111,0,0,5),
    @('San Martin (Parte Holandesa)',549,0,430,100,0,0,19),
    @('Papua Nueva Guinea',511,0,232,273,0,0,6),
    @('Tanzania',509,0,183,305,0,0,21),
    @('Taiwan',499,0,476,16,0,0,7),
    @('Burundi',472,0,374,97,0,0,1),
    @('Comoras',457,0,427,23,0,0,7),
    @('Islas Feroe',423,0,412,11,0,0,0),
    @('Eritrea',361,0,304,57,0,0,0),
    @('Mauricio',361,0,335,16,0,0,10),
    @('Isla de Man',339,0,312,3,0,0,24),
    @('Gibraltar',334,4,307,27,0,0,0),
    @('San Martin (Parte Francesa)',330,0,206,118,0,0,6),
    @('Mongolia',311,0,301,10,0,0,0),
    @('Camboya',275,0,274,1,0,0,0),
    @('Butan',246,1,173,73,0,0,0),
    @('Islas Caimanes',208,0,204,3,0,0,1),
    @('Barbados',183,0,165,11,0,0,7),
    @('Monaco',177,0,132,44,0,0,1),
    @('Bermudas',177,0,161,7,0,0,9),
    @('Curazao',161,0,58,102,0,0,1),
    @('Brunei',145,0,139,3,0,0,3),
    @('Seychelles',140,0,136,4,0,0,0),
    @('Liechtenstein',111,0,105,5,0,0,1),
    @('Antigua y Barbuda',95,0,91,1,0,0,3),
    @('Islas Virgenes Britanicas',66,0,37,28,0,0,1),
    @('San Vicente y las Granadinas',64,0,61,3,0,0,0),
    @('Macao',46,0,46,0,0,0,0),
    @('Puerto Rico',39,0,1,36,0,0,2),
    @('Guam',32,0,0,31,0,0,1),
    @('Fiyi',32,0,24,6,0,0,2),
    @('Bonaire, San Eustaquio y Saba',28,3,17,10,0,0,1),
    @('Santa Lucia',27,0,26,1,0,0,0),
    @('Timor Oriental',27,0,26,1,0,0,0),
    @('Nueva Caledonia',26,0,26,0,0,0,0),
    @('Dominica',24,0,18,6,0,0,0),
    @('Granada',24,0,24,0,0,0,0),
    @('San Bartolome',23,0,16,7,0,0,0),
    @('Laos',23,0,22,1,0,0,0),
    @('Islas Virgenes de los Estados Unidos',17,0,0,17,0,0,0),
    @('San Cristobal y Nieves',17,0,17,0,0,0,0),
    @('Groenlandia',14,0,14,0,0,0,0),
    @('Islas Malvinas',13,0,13,0,0,0,0),
    @('Montserrat',13,0,12,0,0,0,1),
    @('Santa Sede',12,0,12,0,0,0,0),
    @('San Pedro y Miquelon',11,0,5,6,0,0,0),
    @('Sahara Occidental',10,0,8,1,0,0,1),
    @('Anguila',3,0,3,0,0,0,0)
)

$startRow = 4
$rowCount = $data.Count
$colCount = 8

$arr = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $rec = $data[$i]
    for ($j = 0; $j -lt $colCount; $j++) {
        $arr[$i,$j] = $rec[$j]
    }
}

$endRow = $startRow + $rowCount - 1
$rng = $ws.Range("A$startRow" + ":H$endRow")
$rng.Value2 = $arr
